# Natmi following Dr Hou advice
# Add an "ECs" (endothelial cells) sending/target cluster to the Il1a-Il1r2 edge table,
# expanding the 3-row result set (M2 sender only) into 8 rows covering both M2 and ECs
# as sending clusters against ECs/FAPs/M2/sCs as target clusters, with refreshed metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "M2"
$ws.Range("B2").Value = "Il1a"
$ws.Range("C2").Value = "Il1r2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 94.03756700000001
$ws.Range("H2").Value = 282.112701
$ws.Range("I2").Value = 0.948690694112009
$ws.Range("J2").Value = 0.948690694112009
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 263.9035463333333
$ws.Range("N2").Value = 791.710639
$ws.Range("O2").Value = 0.9572387917213622
$ws.Range("P2").Value = 0.9572387917213622
$ws.Range("Q2").Value = 24816.84741985844
$ws.Range("R2").Value = 223351.6267787259
$ws.Range("S2").Value = 0.9081235337490799
$ws.Range("T2").Value = 0.9081235337490799

# Row 3
$ws.Range("A3").Value = "M2"
$ws.Range("B3").Value = "Il1a"
$ws.Range("C3").Value = "Il1r2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 94.03756700000001
$ws.Range("H3").Value = 282.112701
$ws.Range("I3").Value = 0.948690694112009
$ws.Range("J3").Value = 0.948690694112009
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.9720173333333334
$ws.Range("N3").Value = 2.916052
$ws.Range("O3").Value = 0.003525730179150291
$ws.Range("P3").Value = 0.003525730179150291
$ws.Range("Q3").Value = 91.40614510849468
$ws.Range("R3").Value = 822.655305976452
$ws.Range("S3").Value = 0.003344827410909748
$ws.Range("T3").Value = 0.003344827410909748

# Row 4
$ws.Range("A4").Value = "M2"
$ws.Range("B4").Value = "Il1a"
$ws.Range("C4").Value = "Il1r2"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 94.03756700000001
$ws.Range("H4").Value = 282.112701
$ws.Range("I4").Value = 0.948690694112009
$ws.Range("J4").Value = 0.948690694112009
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 9.990833333333333
$ws.Range("N4").Value = 29.9725
$ws.Range("O4").Value = 0.03623904779290016
$ws.Range("P4").Value = 0.03623904779290017
$ws.Range("Q4").Value = 939.5136589691667
$ws.Range("R4").Value = 8455.6229307225
$ws.Range("S4").Value = 0.03437964740460472
$ws.Range("T4").Value = 0.03437964740460472

# Row 5
$ws.Range("A5").Value = "M2"
$ws.Range("B5").Value = "Il1a"
$ws.Range("C5").Value = "Il1r2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 94.03756700000001
$ws.Range("H5").Value = 282.112701
$ws.Range("I5").Value = 0.948690694112009
$ws.Range("J5").Value = 0.948690694112009
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8260933333333332
$ws.Range("N5").Value = 2.47828
$ws.Range("O5").Value = 0.002996430306587325
$ws.Range("P5").Value = 0.002996430306587325
$ws.Range("Q5").Value = 77.68380718158667
$ws.Range("R5").Value = 699.15426463428
$ws.Range("S5").Value = 0.00284268554741459
$ws.Range("T5").Value = 0.00284268554741459

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Il1a"
$ws.Range("C6").Value = "Il1r2"
$ws.Range("D6").Value = "M2"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.085959333333333
$ws.Range("H6").Value = 15.257878
$ws.Range("I6").Value = 0.05130930588799102
$ws.Range("J6").Value = 0.05130930588799102
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 263.9035463333333
$ws.Range("N6").Value = 791.710639
$ws.Range("O6").Value = 0.9572387917213622
$ws.Range("P6").Value = 0.9572387917213622
$ws.Range("Q6").Value = 1342.202704573783
$ws.Range("R6").Value = 12079.82434116404
$ws.Range("S6").Value = 0.0491152579722823
$ws.Range("T6").Value = 0.0491152579722823

# Row 7
$ws.Range("A7").Value = "ECs"
$ws.Range("B7").Value = "Il1a"
$ws.Range("C7").Value = "Il1r2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.085959333333333
$ws.Range("H7").Value = 15.257878
$ws.Range("I7").Value = 0.05130930588799102
$ws.Range("J7").Value = 0.05130930588799102
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.9720173333333334
$ws.Range("N7").Value = 2.916052
$ws.Range("O7").Value = 0.003525730179150291
$ws.Range("P7").Value = 0.003525730179150291
$ws.Range("Q7").Value = 4.943640628628445
$ws.Range("R7").Value = 44.492765657656
$ws.Range("S7").Value = 0.0001809027682405437
$ws.Range("T7").Value = 0.0001809027682405437

# Row 8
$ws.Range("A8").Value = "ECs"
$ws.Range("B8").Value = "Il1a"
$ws.Range("C8").Value = "Il1r2"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.085959333333333
$ws.Range("H8").Value = 15.257878
$ws.Range("I8").Value = 0.05130930588799102
$ws.Range("J8").Value = 0.05130930588799102
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 9.990833333333333
$ws.Range("N8").Value = 29.9725
$ws.Range("O8").Value = 0.03623904779290016
$ws.Range("P8").Value = 0.03623904779290017
$ws.Range("Q8").Value = 50.81297203944444
$ws.Range("R8").Value = 457.316748355
$ws.Range("S8").Value = 0.00185940038829544
$ws.Range("T8").Value = 0.00185940038829544

# Row 9
$ws.Range("A9").Value = "ECs"
$ws.Range("B9").Value = "Il1a"
$ws.Range("C9").Value = "Il1r2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.085959333333333
$ws.Range("H9").Value = 15.257878
$ws.Range("I9").Value = 0.05130930588799102
$ws.Range("J9").Value = 0.05130930588799102
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.8260933333333332
$ws.Range("N9").Value = 2.47828
$ws.Range("O9").Value = 0.002996430306587325
$ws.Range("P9").Value = 0.002996430306587325
$ws.Range("Q9").Value = 4.201477098871111
$ws.Range("R9").Value = 37.81329388984
$ws.Range("S9").Value = 0.0001537447591727358
$ws.Range("T9").Value = 0.0001537447591727358
